$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.144.66"
$ws.Range("D3").Value = "1.577.84"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'209.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "'0.497"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.21%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "'0.0609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").Value = "'19.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "'0.0843"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.799.73"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "1.573.30"
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "'64.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "26.150.12"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "'7.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Value = "'207.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").Value = "'2.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "'143.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "'15.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "'3.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "1.278.18"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'0.609"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.37%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0166"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.93%  "
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D43").Value = "'0.763"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").Value = "'62.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").Value = "1.713.19"
$ws.Range("D46").Value = "'88.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "'1.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").Value = "'0.0506"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("E51").Value = "  -0.18%  "
